$d = $word.ActiveDocument

# Useful characters that are safer to reference by code point than as
# literal source characters (avoids any encoding round-trip issues).
$lq   = [char]8220   # “
$rq   = [char]8221   # ”
$dash = [char]8211   # – (en dash)

# ---------------------------------------------------------------------
# 1) The 17/10/2023 row's "Attivita" cell had its text split across two
#    runs with a "_GoBack" bookmark sitting in between (left over from
#    the previous save). Re-doing the same replacement merges it back
#    into a single run and drops the stale bookmark.
# ---------------------------------------------------------------------
$oldText = "generatore_gif_lst_campania.js" + $rq
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $oldText, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Append the new activity rows right after the 17/10/2023 row (i.e.
#    before the row that used to be the first trailing blank row).
# ---------------------------------------------------------------------
$t = $d.Tables.Item(1)
$insertBeforeIndex = 14

$rowsData = @(
  @("18/10/2023", "10:30 " + $dash + " 16:30", "6", "Revisione codice"),
  @("19/10/2023", "11:00 " + $dash + " 16:00", "5", "Revisione codice e output"),
  @("20/10/2023", "12:00 " + $dash + " 16:00", "4", "Lettura documentazione"),
  @("27/10/2023", "10:30 " + $dash + " 16:30", "6", "Revisione codice e output"),
  @("30/10/2023", "10:30 " + $dash + " 16:30", "6", "Creazione dello script " + $lq + "ndvi_campania_classificato.js" + $rq),
  @("06/11/2023", "10:30 " + $dash + " 16:30", "6", ""),
  @("07/11/2023", "10:30 " + $dash + " 16:30", "6", "")
)

for ($i = $rowsData.Length - 1; $i -ge 0; $i--) {
  $rowData = $rowsData[$i]
  $beforeRow = $t.Rows.Item($insertBeforeIndex)
  $t.Rows.Add($beforeRow) | Out-Null
  $row = $t.Rows.Item($insertBeforeIndex)
  $row.Cells.Item(1).Range.Text = $rowData[0]
  $row.Cells.Item(2).Range.Text = $rowData[1]
  $row.Cells.Item(3).Range.Text = $rowData[2]
  if ($rowData[3] -ne "") {
    $row.Cells.Item(4).Range.Text = $rowData[3]
  }
}

# ---------------------------------------------------------------------
# 3) Update the hour total, then re-plant the "_GoBack" bookmark right
#    after the new text (this mirrors where Word leaves it after the
#    last edit). A zero-length bookmark can't be dropped directly in
#    the slot right before a paragraph mark, so a one-character
#    "shield" is inserted, the bookmark is added next to it, and the
#    shield is deleted again - the bookmark stays put.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Totale ore: 70", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "Totale ore: 109", 2) | Out-Null

# NB: deliberately not using $d.Paragraphs.Last here - after the row
# inserts above, that collection's indexing gets stale in this runtime.
# Find.Execute on a fresh Range keeps working correctly, so use that to
# re-locate the (now updated) total-hours text and anchor off its range.
$totalsRange = $d.Content
$totalsRange.Find.Execute("Totale ore: 109") | Out-Null
$posEnd = $totalsRange.End

$shieldRange = $d.Range($posEnd, $posEnd)
$shieldRange.InsertAfter("X")

$bmRange = $d.Range($posEnd, $posEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$d.Range($posEnd, $posEnd + 1).Delete() | Out-Null
